# Week 15 simulations: append new simulation run numbers to the long
# space-separated number strings on the YDS and ST sheets, and bump the
# aggregate totals on OFF, DEF, ST, TURNS and PEN accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet — append new numbers to the 4 long simulation strings
# ---------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value2 = $wsYDS.Range("B2").Value2 + " 2 3 -2 7 2 -1 3 -2 9 3 4 4 23 8 3 11 13 1 1 2 3 6 -2"
$wsYDS.Range("C2").Value2 = $wsYDS.Range("C2").Value2 + " 20 4 30 14 3 3 16 15 3 29 9 3 0 1 7 7 9 8 5 -4 7 14 -3 5 -2 8 5 4 1 11 2 4 -1 0 2"
$wsYDS.Range("B3").Value2 = $wsYDS.Range("B3").Value2 + " 7 6 15 6 5 6 6 17 29 5 -2 8 3 14 1 37 30 6 15 6 38 9 1 13 3 10 12"
$wsYDS.Range("C3").Value2 = $wsYDS.Range("C3").Value2 + " 8 5 8 14 18 22 19 17 13 0 1 62 12 17"

# ---------------------------------------------------------------------
# OFF sheet — bump aggregate totals
# ---------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value2 = 309
$wsOFF.Range("D2").Value2 = 26
$wsOFF.Range("E2").Value2 = 22
$wsOFF.Range("F2").Value2 = 108
$wsOFF.Range("G2").Value2 = 81
$wsOFF.Range("J2").Value2 = 48
$wsOFF.Range("N2").Value2 = 25
$wsOFF.Range("O2").Value2 = 36
$wsOFF.Range("P2").Value2 = 19

$wsOFF.Range("B3").Value2 = 16
$wsOFF.Range("C3").Value2 = 395
$wsOFF.Range("E3").Value2 = 55
$wsOFF.Range("F3").Value2 = 231
$wsOFF.Range("G3").Value2 = 94
$wsOFF.Range("H3").Value2 = 57
$wsOFF.Range("I3").Value2 = 135
$wsOFF.Range("J3").Value2 = 141
$wsOFF.Range("L3").Value2 = 555
$wsOFF.Range("M3").Value2 = 350
$wsOFF.Range("Q3").Value2 = 882

# ---------------------------------------------------------------------
# DEF sheet — bump aggregate totals
# ---------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value2 = 388
$wsDEF.Range("D2").Value2 = 26
$wsDEF.Range("F2").Value2 = 110
$wsDEF.Range("G2").Value2 = 85
$wsDEF.Range("I2").Value2 = 17
$wsDEF.Range("J2").Value2 = 53
$wsDEF.Range("O2").Value2 = 36
$wsDEF.Range("P2").Value2 = 14

$wsDEF.Range("B3").Value2 = 19
$wsDEF.Range("C3").Value2 = 284
$wsDEF.Range("E3").Value2 = 72
$wsDEF.Range("F3").Value2 = 191
$wsDEF.Range("G3").Value2 = 57
$wsDEF.Range("H3").Value2 = 69
$wsDEF.Range("I3").Value2 = 105
$wsDEF.Range("J3").Value2 = 92
$wsDEF.Range("L3").Value2 = 471
$wsDEF.Range("M3").Value2 = 276
$wsDEF.Range("Q3").Value2 = 915

# ---------------------------------------------------------------------
# ST sheet — append new numbers to the 4 long simulation strings and
# bump the aggregate totals
# ---------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value2 = 158
$wsST.Range("D2").Value2 = 136
$wsST.Range("F2").Value2 = 240
$wsST.Range("G2").Value2 = 226
$wsST.Range("H2").Value2 = 13
$wsST.Range("I2").Value2 = 5
$wsST.Range("L2").Value2 = 71

$wsST.Range("B3").Value2 = 80

$wsST.Range("D3").Value2 = $wsST.Range("D3").Value2 + " 40 49 45 49 41 34"
$wsST.Range("D4").Value2 = $wsST.Range("D4").Value2 + " 16 12 0 0 15 0"
$wsST.Range("D5").Value2 = $wsST.Range("D5").Value2 + " 0 8 0"
$wsST.Range("B6").Value2 = $wsST.Range("B6").Value2 + " 23 33 17"

# ---------------------------------------------------------------------
# TURNS sheet — bump aggregate totals
# ---------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B3").Value2 = 13
$wsTURNS.Range("C3").Value2 = 14
$wsTURNS.Range("D3").Value2 = 14
$wsTURNS.Range("E3").Value2 = 12

# ---------------------------------------------------------------------
# PEN sheet — bump aggregate totals
# ---------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value2 = 26
$wsPEN.Range("B3").Value2 = 31
